# Edit script: apply the changes described in the diff to before.xlsx
# 1) Add new "xml" commands: beautify(xml,var), minify(xml,var) -> sheet "#system" column AA
# 2) Add new "web" command: dragTo(fromLocator,xOffset,yOffset) -> sheet "#system" column V
# 3) Update defined names "web" and "xml" ranges to reflect the new row counts

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --- "web" named range (column V): insert dragTo(fromLocator,xOffset,yOffset) right after
# dragAndDrop(fromLocator,toLocator) at V62, keeping the list alphabetically sorted. ---
# Shift existing V63:V119 down to V64:V120 (processed bottom-up to avoid overwriting).
$ws.Range("V120").Value = "waitForTitle(text)"
$ws.Range("V119").Value = "waitForTextPresent(text)"
$ws.Range("V118").Value = "waitForPopUp(winId,waitMs)"
$ws.Range("V117").Value = "waitForElementPresent(locator)"
$ws.Range("V116").Value = "wait(waitMs)"
$ws.Range("V115").Value = "verifyText(locator,text)"
$ws.Range("V114").Value = "verifyContainText(locator,text)"
$ws.Range("V113").Value = "upload(fieldLocator,file)"
$ws.Range("V112").Value = "unselectAllText()"
$ws.Range("V111").Value = "uncheckAll(locator)"
$ws.Range("V110").Value = "typeKeys(locator,value)"
$ws.Range("V109").Value = "type(locator,value)"
$ws.Range("V108").Value = "toggleSelections(locator)"
$ws.Range("V107").Value = "selectWindowByIndexAndWait(index,waitMs)"
$ws.Range("V106").Value = "selectWindowByIndex(index)"
$ws.Range("V105").Value = "selectWindowAndWait(winId,waitMs)"
$ws.Range("V104").Value = "selectWindow(winId)"
$ws.Range("V103").Value = "selectText(locator)"
$ws.Range("V102").Value = "selectMultiOptions(locator)"
$ws.Range("V101").Value = "selectMulti(locator,array)"
$ws.Range("V100").Value = "selectFrame(locator)"
$ws.Range("V99").Value = "select(locator,text)"
$ws.Range("V98").Value = "scrollTo(locator)"
$ws.Range("V97").Value = "scrollRight(locator,pixel)"
$ws.Range("V96").Value = "scrollLeft(locator,pixel)"
$ws.Range("V95").Value = "saveValue(var,locator)"
$ws.Range("V94").Value = "saveTextSubstringBetween(var,locator,start,end)"
$ws.Range("V93").Value = "saveTextSubstringBefore(var,locator,delim)"
$ws.Range("V92").Value = "saveTextSubstringAfter(var,locator,delim)"
$ws.Range("V91").Value = "saveTextArray(var,locator)"
$ws.Range("V90").Value = "saveText(var,locator)"
$ws.Range("V89").Value = "saveTableAsCsv(locator,nextPageLocator,file)"
$ws.Range("V88").Value = "savePageAsFile(sessionIdName,url,file)"
$ws.Range("V87").Value = "savePageAs(var,sessionIdName,url)"
$ws.Range("V86").Value = "saveLocation(var)"
$ws.Range("V85").Value = "saveLocalStorage(var,key)"
$ws.Range("V84").Value = "saveElements(var,locator)"
$ws.Range("V83").Value = "saveElement(var,locator)"
$ws.Range("V82").Value = "saveDivsAsCsv(headers,rows,cells,nextPage,file)"
$ws.Range("V81").Value = "saveCount(var,locator)"
$ws.Range("V80").Value = "saveAttribute(var,locator,attrName)"
$ws.Range("V79").Value = "saveAllWindowNames(var)"
$ws.Range("V78").Value = "saveAllWindowIds(var)"
$ws.Range("V77").Value = "resizeWindow(width,height)"
$ws.Range("V76").Value = "refreshAndWait()"
$ws.Range("V75").Value = "refresh()"
$ws.Range("V74").Value = "openIgnoreTimeout(url)"
$ws.Range("V73").Value = "openHttpBasic(url,username,password)"
$ws.Range("V72").Value = "openAndWait(url,waitMs)"
$ws.Range("V71").Value = "open(url)"
$ws.Range("V70").Value = "mouseOver(locator)"
$ws.Range("V69").Value = "maximizeWindow()"
$ws.Range("V68").Value = "goBackAndWait()"
$ws.Range("V67").Value = "goBack()"
$ws.Range("V66").Value = "focus(locator)"
$ws.Range("V65").Value = "executeScript(var,script)"
$ws.Range("V64").Value = "editLocalStorage(key,value)"

# Insert the new command into the now-vacant V63 slot.
$ws.Range("V63").Value = "dragTo(fromLocator,xOffset,yOffset)"

# --- "xml" named range (column AA): insert beautify(xml,var) and minify(xml,var) right after
# assertWellformed(xml) at AA8, keeping the list alphabetically sorted. ---
# Shift existing AA9:AA11 down to AA11:AA13 (processed bottom-up to avoid overwriting).
$ws.Range("AA13").Value = "storeValues(xml,xpath,var)"
$ws.Range("AA12").Value = "storeValue(xml,xpath,var)"
$ws.Range("AA11").Value = "storeCount(xml,xpath,var)"

# Insert the two new commands into the now-vacant AA9 and AA10 slots.
$ws.Range("AA9").Value = "beautify(xml,var)"
$ws.Range("AA10").Value = "minify(xml,var)"

# --- Update the defined names so they cover the newly extended lists. ---
$wb.Names.Item("web").RefersTo = "='#system'!`$V`$2:`$V`$120"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AA`$2:`$AA`$13"

